$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the casing of a few label strings in column A ("classes" names)
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A8").Value = "pageTitleNewTab"
$ws.Range("A4").Value = "mdaTitle"

# Update the active selection to A4
$ws.Range("A4").Select()
